# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows that changed when the data
# was repulled from source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2  = 5
    7  = -3
    8  = 1
    9  = -2
    11 = 0
    23 = -2
    24 = 0
    27 = -2
    29 = -4
    30 = -3
    33 = 0
    34 = -4
    37 = 2
    44 = 1
    46 = -1
    47 = -1
    51 = -1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
